$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.922.44"
$ws.Range("E2").Value = "  -0.39%  "

$ws.Range("D3").Value = "3.566.97"
$ws.Range("E3").Value = "  +2.52%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.76%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.25%  "

$ws.Range("D7").Value = "3.566.25"
$ws.Range("E7").Value = "  +2.51%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  +0.51%  "

$ws.Range("E10").Value = "  +0.39%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.97"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.55%  "

$ws.Range("E12").Value = "  -0.01%  "

$ws.Range("D13").Value = "4.171.98"
$ws.Range("E13").Value = "  +2.52%  "

$ws.Range("E14").Value = "  -0.04%  "

$ws.Range("D15").Value = "3.570.78"
$ws.Range("E15").Value = "  +3.45%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.30%  "

$ws.Range("E17").Value = "  +0.46%  "

$ws.Range("D18").Value = "65.058.78"
$ws.Range("E18").Value = "  -0.05%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.71%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.84%  "

$ws.Range("E21").Value = "  +0.61%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "388.91"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.15%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.578"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.02%  "

$ws.Range("D24").Value = "3.711.64"
$ws.Range("E24").Value = "  +2.55%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.83%  "

$ws.Range("E26").Value = "  +0.21%  "

$ws.Range("E27").Value = "  +3.96%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.89%  "

$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.40"
$ws.Range("D31").Style = "Normal"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.49"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +24.06%  "

$ws.Range("D33").Value = "3.567.65"
$ws.Range("E33").Value = "  +1.98%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.02"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.96%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.144"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.27%  "

$ws.Range("E37").Value = "  +1.41%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "168.75"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.45%  "

$ws.Range("E39").Value = "  +4.30%  "

$ws.Range("E40").Value = "  +5.37%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0805"
$ws.Range("D41").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.825"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.58%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.84"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.36%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.29%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.01%  "

$ws.Range("E46").Value = "  +1.98%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.75%  "

$ws.Range("E48").Value = "  +1.49%  "

$ws.Range("D49").Value = "2.481.83"
$ws.Range("E49").Value = "  +11.95%  "

$ws.Range("E50").Value = "  +2.77%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.865"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.68%  "
